$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: beetle sample BE64 (row 4) was re-measured ---
$ws.Range("B4").Value = 20.7

# --- Column C: "Long"/"Short" result formula for every sample ---
# C4 is filled individually (matches the diff, which shows C4 as a standalone
# formula while C5:C20 become one shared-formula group).
$ws.Range("C4").Formula = '=IF($B4>$B$21,"Long","Short")'
$ws.Range("C5:C20").Formula = '=IF($B5>$B$21,"Long","Short")'

# --- Column D: extra comment formula when a sample is "Long" ---
$ws.Range("D4").Formula = '=IF(C4="Long","This sample is "&ROUND($B4-$B$21,2)&" longer than avrg","")'
$ws.Range("D5:D20").Formula = '=IF(C5="Long","This sample is "&ROUND($B5-$B$21,2)&" longer than avrg","")'

# --- Column D is widened to fit the longer comment text ---
$ws.Columns("D").ColumnWidth = 27

# --- View state: zoomed in, selection moved ---
$ws.Application.ActiveWindow.Zoom = 140
[void]$ws.Range("H18").Select()
